# Apply table formatting changes to the first table in the document:
#  - every cell gets a thin light-gray (D3D3D3) single border on all 4 sides
#  - header row (row 1) cell shading changes from FFEFD5 to FFDAB9
#  - data rows (rows 2-9), columns 2-4 paragraph alignment changes from
#    right-justified to centered

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$wdBorderTop = -1
$wdBorderLeft = -2
$wdBorderBottom = -3
$wdBorderRight = -4

$wdLineStyleSingle = 1
$wdLineWidth025pt = 2

$borderColor = 0xD3D3D3          # BGR == RGB here since the color is gray
$headerFill  = 0xB9DAFF          # BGR encoding of target RGB FFDAB9
$wdAlignParagraphCenter = 1

$rowCount = $t.Rows.Count
$colCount = $t.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $t.Cell($r, $c)

        $cell.Borders.Item($wdBorderTop).LineStyle = $wdLineStyleSingle
        $cell.Borders.Item($wdBorderTop).LineWidth = $wdLineWidth025pt
        $cell.Borders.Item($wdBorderTop).Color = $borderColor

        $cell.Borders.Item($wdBorderLeft).LineStyle = $wdLineStyleSingle
        $cell.Borders.Item($wdBorderLeft).LineWidth = $wdLineWidth025pt
        $cell.Borders.Item($wdBorderLeft).Color = $borderColor

        $cell.Borders.Item($wdBorderBottom).LineStyle = $wdLineStyleSingle
        $cell.Borders.Item($wdBorderBottom).LineWidth = $wdLineWidth025pt
        $cell.Borders.Item($wdBorderBottom).Color = $borderColor

        $cell.Borders.Item($wdBorderRight).LineStyle = $wdLineStyleSingle
        $cell.Borders.Item($wdBorderRight).LineWidth = $wdLineWidth025pt
        $cell.Borders.Item($wdBorderRight).Color = $borderColor

        if ($r -eq 1) {
            $cell.Shading.BackgroundPatternColor = $headerFill
        } elseif ($c -ge 2) {
            $cell.Range.ParagraphFormat.Alignment = $wdAlignParagraphCenter
        }
    }
}

Write-Output "table formatting updated"
